$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''28.826.57'
$ws.Range("E2").Value = '''  +2.49%  '

# Row 3
$ws.Range("D3").Value = '''1.881.39'
$ws.Range("E3").Value = '''  +2.70%  '

# Row 4
$ws.Range("E4").Value = '''  +0.60%  '

# Row 5
$ws.Range("D5").Value = '''324.29'
$ws.Range("E5").Value = '''  -1.35%  '

# Row 6
$ws.Range("D6").Value = '''1.006'
$ws.Range("E6").Value = '''  +0.60%  '

# Row 7
$ws.Range("D7").Value = '''0.4675'

# Row 8
$ws.Range("D8").Value = '''0.3933'
$ws.Range("E8").Value = '''  +1.65%  '

# Row 9
$ws.Range("D9").Value = '''0.07928'

# Row 10
$ws.Range("D10").Value = '''0.9818'
$ws.Range("E10").Value = '''  +2.13%  '

# Row 11
$ws.Range("D11").Value = '''22.39'
$ws.Range("E11").Value = '''  +1.89%  '

# Row 12
$ws.Range("D12").Value = '''1.836.42'
$ws.Range("E12").Value = '''  +1.92%  '

# Row 13
$ws.Range("D13").Value = '''5.744'
$ws.Range("E13").Value = '''  +1.33%  '

# Row 14
$ws.Range("D14").Value = '''7.018'
$ws.Range("E14").Value = '''  +1.72%  '

# Row 15
$ws.Range("D15").Value = '''0.06982'
$ws.Range("E15").Value = '''  +1.88%  '

# Row 16
$ws.Range("D16").Value = '''88.77'
$ws.Range("E16").Value = '''  +2.44%  '

# Row 17
$ws.Range("E17").Value = '''  +0.66%  '

# Row 18
$ws.Range("D18").Value = '''0.00001011'
$ws.Range("E18").Value = '''  +1.48%  '

# Row 19
$ws.Range("D19").Value = '''16.98'
$ws.Range("E19").Value = '''  +1.86%  '

# Row 20
$ws.Range("E20").Value = '''  +0.46%  '

# Row 21
$ws.Range("D21").Value = '''28.853.18'
$ws.Range("E21").Value = '''  +2.56%  '

# Row 22
$ws.Range("D22").Value = '''5.351'
$ws.Range("E22").Value = '''  +0.33%  '

# Row 23
$ws.Range("E23").Value = '''  +0.74%  '

# Row 24
$ws.Range("B24").Value = '''WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = '''https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '''2.188.63'
$ws.Range("E24").Value = '''  +7.54%  '

# Row 25
$ws.Range("B25").Value = '''Toncoin'
$ws.Range("C25").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''2.129'
$ws.Range("E25").Value = '''  +1.74%  '

# Row 26
$ws.Range("D26").Value = '''153.38'
$ws.Range("E26").Value = '''  +0.63%  '

# Row 27
$ws.Range("D27").Value = '''19.42'
$ws.Range("E27").Value = '''  +1.10%  '

# Row 28
$ws.Range("D28").Value = '''5.804'
$ws.Range("E28").Value = '''  +0.33%  '

# Row 29
$ws.Range("D29").Value = '''2.002'
$ws.Range("E29").Value = '''  +1.38%  '

# Row 30
$ws.Range("D30").Value = '''119.93'
$ws.Range("E30").Value = '''  +2.51%  '

# Row 31
$ws.Range("D31").Value = '''0.09399'
$ws.Range("E31").Value = '''  +1.51%  '

# Row 32
$ws.Range("D32").Value = '''0.9378'
$ws.Range("E32").Value = '''  +0.08%  '

# Row 33
$ws.Range("D33").Value = '''5.319'
$ws.Range("E33").Value = '''  +0.24%  '

# Row 34
$ws.Range("D34").Value = '''1.358'
$ws.Range("E34").Value = '''  +2.84%  '

# Row 35
$ws.Range("D35").Value = '''3.350'
$ws.Range("E35").Value = '''  +0.03%  '

# Row 36
$ws.Range("D36").Value = '''0.05914'
$ws.Range("E36").Value = '''  -0.52%  '

# Row 37
$ws.Range("D37").Value = '''0.02130'
$ws.Range("E37").Value = '''  -0.91%  '

# Row 38
$ws.Range("D38").Value = '''1.163'
$ws.Range("E38").Value = '''  +1.15%  '

# Row 39
$ws.Range("D39").Value = '''7.900'
$ws.Range("E39").Value = '''  +3.23%  '

# Row 40
$ws.Range("D40").Value = '''0.5733'
$ws.Range("E40").Value = '''  +2.53%  '

# Row 41
$ws.Range("B41").Value = '''Aptos'
$ws.Range("C41").Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '''10.01'
$ws.Range("E41").Value = '''  +0.77%  '

# Row 42
$ws.Range("B42").Value = '''Algorand'
$ws.Range("C42").Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.1798'
$ws.Range("E42").Value = '''  +1.34%  '

# Row 43
$ws.Range("D43").Value = '''0.07316'
$ws.Range("E43").Value = '''  +3.85%  '

# Row 44
$ws.Range("D44").Value = '''11.82'
$ws.Range("E44").Value = '''  +1.79%  '

# Row 45
$ws.Range("D45").Value = '''1.180'
$ws.Range("E45").Value = '''  -3.47%  '

# Row 46
$ws.Range("D46").Value = '''0.5361'
$ws.Range("E46").Value = '''  +1.55%  '

# Row 47
$ws.Range("D47").Value = '''1.846'
$ws.Range("E47").Value = '''  +0.60%  '

# Row 48
$ws.Range("D48").Value = '''114.03'
$ws.Range("E48").Value = '''  +2.08%  '

# Row 49
$ws.Range("D49").Value = '''2.074'
$ws.Range("E49").Value = '''  -6.48%  '

# Row 50
$ws.Range("D50").Value = '''2.375'
$ws.Range("E50").Value = '''  +2.76%  '

# Row 51
$ws.Range("D51").Value = '''1.005'
$ws.Range("E51").Value = '''  +0.60%  '
